$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.346.30"
$ws.Range("E2").Value = "  -2.38%  "

$ws.Range("D3").Value = "1.709.30"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5313"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2657"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.571"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.72%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.948.87"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.660.03"
$ws.Range("E14").Value = "  -4.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5720"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.13%  "

$ws.Range("D16").Value = "0.0₅8165"
$ws.Range("E16").Value = "  -2.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").Value = "27.373.06"
$ws.Range("E18").Value = "  -2.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.671"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.970"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.768"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1217"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.274"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05413"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.293"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.502"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.427"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.643"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.882"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9490"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.417"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D41").Value = "1.046.32"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8449"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").Value = "1.854.55"

$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4506"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.091"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05246"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
